$d = $word.ActiveDocument

$null = $d.Content.Find.Execute("26+38=64", $true, $true, $false, $false, $false, $true, 1, $false, "5+62=67", 2)
$null = $d.Content.Find.Execute("84-28=56", $true, $true, $false, $false, $false, $true, 1, $false, "31+38=69", 2)
$null = $d.Content.Find.Execute("84-50=34", $true, $true, $false, $false, $false, $true, 1, $false, "47-36=11", 2)
$null = $d.Content.Find.Execute("64-9=55", $true, $true, $false, $false, $false, $true, 1, $false, "62+16=78", 2)
$null = $d.Content.Find.Execute("82-79=3", $true, $true, $false, $false, $false, $true, 1, $false, "63+19=82", 2)
$null = $d.Content.Find.Execute("26-24=2", $true, $true, $false, $false, $false, $true, 1, $false, "16+41=57", 2)
$null = $d.Content.Find.Execute("51+7=58", $true, $true, $false, $false, $false, $true, 1, $false, "55-6=49", 2)
$null = $d.Content.Find.Execute("70-9=61", $true, $true, $false, $false, $false, $true, 1, $false, "64+9=73", 2)
$null = $d.Content.Find.Execute("59+6=65", $true, $true, $false, $false, $false, $true, 1, $false, "61-24=37", 2)
$null = $d.Content.Find.Execute("13+61=74", $true, $true, $false, $false, $false, $true, 1, $false, "59-2=57", 2)
$null = $d.Content.Find.Execute("64-21=43", $true, $true, $false, $false, $false, $true, 1, $false, "34+33=67", 2)
$null = $d.Content.Find.Execute("24+71=95", $true, $true, $false, $false, $false, $true, 1, $false, "42-3=39", 2)
$null = $d.Content.Find.Execute("21-1=20", $true, $true, $false, $false, $false, $true, 1, $false, "33+57=90", 2)
$null = $d.Content.Find.Execute("53+4=57", $true, $true, $false, $false, $false, $true, 1, $false, "94-9=85", 2)
$null = $d.Content.Find.Execute("43+34=77", $true, $true, $false, $false, $false, $true, 1, $false, "16+68=84", 2)
$null = $d.Content.Find.Execute("56-9=47", $true, $true, $false, $false, $false, $true, 1, $false, "79-52=27", 2)
$null = $d.Content.Find.Execute("1+55=56", $true, $true, $false, $false, $false, $true, 1, $false, "18-6=12", 2)
$null = $d.Content.Find.Execute("97-55=42", $true, $true, $false, $false, $false, $true, 1, $false, "86-42=44", 2)
$null = $d.Content.Find.Execute("30+32=62", $true, $true, $false, $false, $false, $true, 1, $false, "18+33=51", 2)
$null = $d.Content.Find.Execute("69-7=62", $true, $true, $false, $false, $false, $true, 1, $false, "54+11=65", 2)
$null = $d.Content.Find.Execute("63+9=72", $true, $true, $false, $false, $false, $true, 1, $false, "39+13=52", 2)
$null = $d.Content.Find.Execute("72-57=15", $true, $true, $false, $false, $false, $true, 1, $false, "30+58=88", 2)
$null = $d.Content.Find.Execute("44-2=42", $true, $true, $false, $false, $false, $true, 1, $false, "55-46=9", 2)
$null = $d.Content.Find.Execute("21-7=14", $true, $true, $false, $false, $false, $true, 1, $false, "76-2=74", 2)
$null = $d.Content.Find.Execute("9-9=0", $true, $true, $false, $false, $false, $true, 1, $false, "16+2=18", 2)
$null = $d.Content.Find.Execute("10+47=57", $true, $true, $false, $false, $false, $true, 1, $false, "74-71=3", 2)
$null = $d.Content.Find.Execute("82-50=32", $true, $true, $false, $false, $false, $true, 1, $false, "25-23=2", 2)
$null = $d.Content.Find.Execute("92-41=51", $true, $true, $false, $false, $false, $true, 1, $false, "54-29=25", 2)
$null = $d.Content.Find.Execute("63-62=1", $true, $true, $false, $false, $false, $true, 1, $false, "87+8=95", 2)
$null = $d.Content.Find.Execute("21+35=56", $true, $true, $false, $false, $false, $true, 1, $false, "65-52=13", 2)
$null = $d.Content.Find.Execute("98-83=15", $true, $true, $false, $false, $false, $true, 1, $false, "30+6=36", 2)
$null = $d.Content.Find.Execute("96-71=25", $true, $true, $false, $false, $false, $true, 1, $false, "11+78=89", 2)
$null = $d.Content.Find.Execute("28+57=85", $true, $true, $false, $false, $false, $true, 1, $false, "20-6=14", 2)
$null = $d.Content.Find.Execute("12+63=75", $true, $true, $false, $false, $false, $true, 1, $false, "88-66=22", 2)
$null = $d.Content.Find.Execute("17+67=84", $true, $true, $false, $false, $false, $true, 1, $false, "95-45=50", 2)
$null = $d.Content.Find.Execute("62+0=62", $true, $true, $false, $false, $false, $true, 1, $false, "16+9=25", 2)
$null = $d.Content.Find.Execute("53+2=55", $true, $true, $false, $false, $false, $true, 1, $false, "51-50=1", 2)
$null = $d.Content.Find.Execute("35+13=48", $true, $true, $false, $false, $false, $true, 1, $false, "27+42=69", 2)
$null = $d.Content.Find.Execute("20+6=26", $true, $true, $false, $false, $false, $true, 1, $false, "85-48=37", 2)
$null = $d.Content.Find.Execute("9+72=81", $true, $true, $false, $false, $false, $true, 1, $false, "39+29=68", 2)
$null = $d.Content.Find.Execute("24-18=6", $true, $true, $false, $false, $false, $true, 1, $false, "14+44=58", 2)
$null = $d.Content.Find.Execute("28+56=84", $true, $true, $false, $false, $false, $true, 1, $false, "43+33=76", 2)
$null = $d.Content.Find.Execute("20+67=87", $true, $true, $false, $false, $false, $true, 1, $false, "51+5=56", 2)
$null = $d.Content.Find.Execute("81-73=8", $true, $true, $false, $false, $false, $true, 1, $false, "20+61=81", 2)
$null = $d.Content.Find.Execute("9+11=20", $true, $true, $false, $false, $false, $true, 1, $false, "82-69=13", 2)
$null = $d.Content.Find.Execute("64-48=16", $true, $true, $false, $false, $false, $true, 1, $false, "33-31=2", 2)
$null = $d.Content.Find.Execute("10-7=3", $true, $true, $false, $false, $false, $true, 1, $false, "47+14=61", 2)
$null = $d.Content.Find.Execute("85-29=56", $true, $true, $false, $false, $false, $true, 1, $false, "90-31=59", 2)
$null = $d.Content.Find.Execute("86-59=27", $true, $true, $false, $false, $false, $true, 1, $false, "92-14=78", 2)
$null = $d.Content.Find.Execute("21+16=37", $true, $true, $false, $false, $false, $true, 1, $false, "26+61=87", 2)
$null = $d.Content.Find.Execute("78-47=31", $true, $true, $false, $false, $false, $true, 1, $false, "57+8=65", 2)
$null = $d.Content.Find.Execute("61-23=38", $true, $true, $false, $false, $false, $true, 1, $false, "79-56=23", 2)
$null = $d.Content.Find.Execute("38-22=16", $true, $true, $false, $false, $false, $true, 1, $false, "63-56=7", 2)
$null = $d.Content.Find.Execute("66+6=72", $true, $true, $false, $false, $false, $true, 1, $false, "72-38=34", 2)
$null = $d.Content.Find.Execute("7+84=91", $true, $true, $false, $false, $false, $true, 1, $false, "76-71=5", 2)
$null = $d.Content.Find.Execute("14+53=67", $true, $true, $false, $false, $false, $true, 1, $false, "48-27=21", 2)
$null = $d.Content.Find.Execute("5+17=22", $true, $true, $false, $false, $false, $true, 1, $false, "11+66=77", 2)
$null = $d.Content.Find.Execute("28-0=28", $true, $true, $false, $false, $false, $true, 1, $false, "20+4=24", 2)
$null = $d.Content.Find.Execute("48+25=73", $true, $true, $false, $false, $false, $true, 1, $false, "18-8=10", 2)
$null = $d.Content.Find.Execute("65-21=44", $true, $true, $false, $false, $false, $true, 1, $false, "51-48=3", 2)
$null = $d.Content.Find.Execute("63-22=41", $true, $true, $false, $false, $false, $true, 1, $false, "43+15=58", 2)
$null = $d.Content.Find.Execute("55+29=84", $true, $true, $false, $false, $false, $true, 1, $false, "85-11=74", 2)
$null = $d.Content.Find.Execute("54+33=87", $true, $true, $false, $false, $false, $true, 1, $false, "96-90=6", 2)
$null = $d.Content.Find.Execute("17+8=25", $true, $true, $false, $false, $false, $true, 1, $false, "27+42=69", 2)
$null = $d.Content.Find.Execute("91-34=57", $true, $true, $false, $false, $false, $true, 1, $false, "15-8=7", 2)
$null = $d.Content.Find.Execute("36-29=7", $true, $true, $false, $false, $false, $true, 1, $false, "10+43=53", 2)
$null = $d.Content.Find.Execute("28+20=48", $true, $true, $false, $false, $false, $true, 1, $false, "67-64=3", 2)
$null = $d.Content.Find.Execute("9+79=88", $true, $true, $false, $false, $false, $true, 1, $false, "10+12=22", 2)
$null = $d.Content.Find.Execute("44+23=67", $true, $true, $false, $false, $false, $true, 1, $false, "91-64=27", 2)
$null = $d.Content.Find.Execute("69+2=71", $true, $true, $false, $false, $false, $true, 1, $false, "49+45=94", 2)
$null = $d.Content.Find.Execute("81-48=33", $true, $true, $false, $false, $false, $true, 1, $false, "97-82=15", 2)
$null = $d.Content.Find.Execute("85-0=85", $true, $true, $false, $false, $false, $true, 1, $false, "35+52=87", 2)
$null = $d.Content.Find.Execute("62-2=60", $true, $true, $false, $false, $false, $true, 1, $false, "62-33=29", 2)
$null = $d.Content.Find.Execute("22+21=43", $true, $true, $false, $false, $false, $true, 1, $false, "90-65=25", 2)
$null = $d.Content.Find.Execute("16+48=64", $true, $true, $false, $false, $false, $true, 1, $false, "58-33=25", 2)
$null = $d.Content.Find.Execute("99-7=92", $true, $true, $false, $false, $false, $true, 1, $false, "4+86=90", 2)
$null = $d.Content.Find.Execute("27+25=52", $true, $true, $false, $false, $false, $true, 1, $false, "34+58=92", 2)
$null = $d.Content.Find.Execute("77+11=88", $true, $true, $false, $false, $false, $true, 1, $false, "90-64=26", 2)
$null = $d.Content.Find.Execute("29+58=87", $true, $true, $false, $false, $false, $true, 1, $false, "67-4=63", 2)
$null = $d.Content.Find.Execute("36+27=63", $true, $true, $false, $false, $false, $true, 1, $false, "41+43=84", 2)
$null = $d.Content.Find.Execute("72-15=57", $true, $true, $false, $false, $false, $true, 1, $false, "2+27=29", 2)
$null = $d.Content.Find.Execute("54+44=98", $true, $true, $false, $false, $false, $true, 1, $false, "56-41=15", 2)
$null = $d.Content.Find.Execute("51-29=22", $true, $true, $false, $false, $false, $true, 1, $false, "64-5=59", 2)
$null = $d.Content.Find.Execute("76-70=6", $true, $true, $false, $false, $false, $true, 1, $false, "81-27=54", 2)
$null = $d.Content.Find.Execute("51+13=64", $true, $true, $false, $false, $false, $true, 1, $false, "32+45=77", 2)
$null = $d.Content.Find.Execute("1+43=44", $true, $true, $false, $false, $false, $true, 1, $false, "65-30=35", 2)
$null = $d.Content.Find.Execute("9+61=70", $true, $true, $false, $false, $false, $true, 1, $false, "74+1=75", 2)
$null = $d.Content.Find.Execute("68+29=97", $true, $true, $false, $false, $false, $true, 1, $false, "73-45=28", 2)
$null = $d.Content.Find.Execute("60-44=16", $true, $true, $false, $false, $false, $true, 1, $false, "1-0=1", 2)
$null = $d.Content.Find.Execute("32-0=32", $true, $true, $false, $false, $false, $true, 1, $false, "11+13=24", 2)
$null = $d.Content.Find.Execute("27-5=22", $true, $true, $false, $false, $false, $true, 1, $false, "81-9=72", 2)
$null = $d.Content.Find.Execute("25+31=56", $true, $true, $false, $false, $false, $true, 1, $false, "77+20=97", 2)
$null = $d.Content.Find.Execute("76-22=54", $true, $true, $false, $false, $false, $true, 1, $false, "21-2=19", 2)
$null = $d.Content.Find.Execute("88-77=11", $true, $true, $false, $false, $false, $true, 1, $false, "80-49=31", 2)
$null = $d.Content.Find.Execute("8+1=9", $true, $true, $false, $false, $false, $true, 1, $false, "71-32=39", 2)
$null = $d.Content.Find.Execute("49-32=17", $true, $true, $false, $false, $false, $true, 1, $false, "80+3=83", 2)
$null = $d.Content.Find.Execute("45+25=70", $true, $true, $false, $false, $false, $true, 1, $false, "32+46=78", 2)
$null = $d.Content.Find.Execute("87-43=44", $true, $true, $false, $false, $false, $true, 1, $false, "47+25=72", 2)
$null = $d.Content.Find.Execute("72+25=97", $true, $true, $false, $false, $false, $true, 1, $false, "71+10=81", 2)
$null = $d.Content.Find.Execute("83-7=76", $true, $true, $false, $false, $false, $true, 1, $false, "79-13=66", 2)
